$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing content first
$ws.Cells.Clear()

$items = @(
    '{"name":"Dreams & Nightmares Case","price":"$1.37 USD"}',
    '{"name":"Revolution Case","price":"$1.43 USD"}',
    '{"name":"Fracture Case","price":"$0.66 USD"}',
    '{"name":"Paris 2023 Legends Sticker Capsule","price":"$0.26 USD"}',
    '{"name":"Mann Co. Supply Crate Key","price":"$2.15 USD"}'
)

$row = 1
foreach ($item in $items) {
    $ws.Cells.Item($row, 1).Value = $item
    $row++
}
